$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("top-interview-questions-easy").Name = "Top Interview Q - Easy 48"
$wb.Worksheets.Item("Meduim Collection 52").Name = "Top Interview Q - Meduim 52"

# ---------------------------------------------------------------------------
# 2) Add the two new sheets at the end: "Graph" first (so it gets the lower
#    sheetId), then "Binary Search" (higher sheetId), then move "Binary
#    Search" so it sits before "Graph" in tab order - matches the sheetId /
#    tab-order split seen in the target workbook.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$graphSheet = $wb.Worksheets.Add($null, $lastSheet)
$graphSheet.Name = "Graph"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$bsSheet = $wb.Worksheets.Add($null, $lastSheet2)
$bsSheet.Name = "Binary Search"

$bsSheet.Move($graphSheet)

# Re-fetch by name after the Move (avoids any stale-handle surprises).
$bsSheet = $wb.Worksheets.Item("Binary Search")
$graphSheet = $wb.Worksheets.Item("Graph")

$bsSheet.Activate()
$bsSheet.Range("R38").Select()

$graphSheet.Activate()
$graphSheet.Range("A2").Select()

# ---------------------------------------------------------------------------
# 3) Update remembered selections on the existing sheets.
# ---------------------------------------------------------------------------
$easySheet = $wb.Worksheets.Item("Top Interview Q - Easy 48")
$easySheet.Activate()
$easySheet.Range("E35").Select()

$blind75 = $wb.Worksheets.Item("Blind 75")
$blind75.Activate()
$blind75.Range("G70").Select()

$mediumSheet = $wb.Worksheets.Item("Top Interview Q - Meduim 52")
$mediumSheet.Activate()
$mediumSheet.Range("F42").Select()

# Restore the workbook's active tab to "Blind 75" (unchanged by this edit).
$blind75.Activate()

# ---------------------------------------------------------------------------
# 4) New Blind 75 rows: #17 Letter Combinations of a Phone Number, and fill
#    in the LC number / notes / date for #11 Container With Most Water.
# ---------------------------------------------------------------------------
$ws = $blind75

# Row 62 - add LC link for "Letter Combinations of a Phone Number" first so
# the new shared string is appended in the same order as the source edit.
$ws.Range("G62").Value = "https://leetcode.com/problems/letter-combinations-of-a-phone-number/"
$ws.Range("G61").Copy()
$ws.Range("G62").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C62").Value = 17

# Row 61 - "Container With Most Water" gets its LC number, solved date and
# approach notes.
$ws.Range("C61").Value = 11

$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = "2023.09.26"
$ws.Range("A61").Copy()
$ws.Range("I61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H61").Value = "two pointer: left, right, find curArea and check maxArea, if left < right height, then left ++ or right--"

$excel.CutCopyMode = $false
